$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")
$ws.Range("B28").Value = "'34"
$ws.Activate()
$ws.Range("B29").Select() | Out-Null
